# Apply row-rotation edits to the Artfynd worksheet.
#
# Each group of row numbers below has its full row contents (columns A:AY)
# rotated by one position: the data that was in the first row of the group
# moves to the second row, the data in the second row moves to the third
# row, and so on, with the data from the last row of the group wrapping
# around to the first row of the group.
#
# Implementation notes:
#  - Each source row is first staged into a scratch row (far below the real
#    data) so that, regardless of write order, a row used as a source for
#    one destination is never clobbered before it has been fully read.
#  - Cells are copied one at a time with Range.Copy (rather than assigning
#    .Value2 from a bulk array) so that each cell keeps its original
#    type/format - notably this avoids Excel's implicit text->date
#    coercion that happens when a date-looking string (e.g. "2026-02-05")
#    is pushed through a plain .Value2 assignment.
#  - A column is only touched (cleared + copied) when either the source or
#    the current destination cell is non-blank; columns that are genuinely
#    absent in both keep the worksheet's original sparse layout instead of
#    materialising empty cell records.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastColIndex = 51 # AY

# Groups of rows whose contents rotate (old row N's data -> row N+1,
# with the last row's data wrapping back around to the first row).
$groups = @(
    , @(2, 3, 4)
    , @(9, 10)
    , @(12, 13)
    , @(16, 17)
    , @(18, 19, 20)
    , @(26, 27)
)

$scratchBase = 1000

foreach ($group in $groups) {
    $n = $group.Count

    # Stage a copy of every row in this group into scratch rows first,
    # so sources remain intact regardless of write order.
    for ($i = 0; $i -lt $n; $i++) {
        $r = $group[$i]
        $scratchRow = $scratchBase + $i
        for ($col = 1; $col -le $lastColIndex; $col++) {
            $srcCell = $ws.Cells.Item($r, $col)
            $scratchCell = $ws.Cells.Item($scratchRow, $col)
            $srcVal = $srcCell.Value2
            $scratchVal = $scratchCell.Value2
            if (($null -eq $srcVal) -and ($null -eq $scratchVal)) {
                continue
            }
            $scratchCell.ClearContents()
            $srcCell.Copy($scratchCell)
        }
    }

    # Write each row's new content: row at index i receives the (staged)
    # snapshot from the row at index i-1 (cyclic).
    for ($i = 0; $i -lt $n; $i++) {
        $destRow = $group[$i]
        $srcIndex = (($i - 1) + $n) % $n
        $scratchRow = $scratchBase + $srcIndex
        for ($col = 1; $col -le $lastColIndex; $col++) {
            $scratchCell = $ws.Cells.Item($scratchRow, $col)
            $destCell = $ws.Cells.Item($destRow, $col)
            $scratchVal = $scratchCell.Value2
            $destVal = $destCell.Value2
            if (($null -eq $scratchVal) -and ($null -eq $destVal)) {
                continue
            }
            $destCell.ClearContents()
            $scratchCell.Copy($destCell)
        }
    }

    # Clean up the scratch rows used for this group.
    for ($i = 0; $i -lt $n; $i++) {
        $scratchRow = $scratchBase + $i
        $scratchRng = $ws.Range("A" + $scratchRow + ":AY" + $scratchRow)
        $scratchRng.Clear()
    }
}
